# Updates the cryptos list (prices / 1h volume %) to match the
# "Updated cryptos list on Wed Mar 20 21:22:31 UTC 2024 with GitHub Actions" commit.
#
# Price/name/link cells are stored as plain text in the workbook (so that
# values such as "67.055.86" or "1.00" are never reinterpreted as numbers).
# Assigning a bare numeric-looking string to Range.Value lets Excel coerce it
# back into a real number (dropping the text formatting), so for any new
# value that Excel would otherwise auto-convert we prefix it with a leading
# apostrophe (a doubled '' inside a PowerShell single-quoted string yields a
# literal '), which is how Excel keeps text entry "as typed" instead of
# reparsing it as a number; the leading marker itself is not stored as part
# of the cell's text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.257.13'
$ws.Range("E2").Value = '  +5.07%  '
$ws.Range("D3").Value = '3.485.78'
$ws.Range("E3").Value = '  +6.15%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''187.01'
$ws.Range("E5").Value = '  +7.39%  '
$ws.Range("D6").Value = '''548.39'
$ws.Range("E6").Value = '  +4.98%  '
$ws.Range("D7").Value = '''0.613'
$ws.Range("E7").Value = '  +1.92%  '
$ws.Range("D8").Value = '3.487.73'
$ws.Range("E8").Value = '  +6.50%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '''0.639'
$ws.Range("E10").Value = '  +5.55%  '
$ws.Range("D11").Value = '''56.45'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '''0.148'
$ws.Range("E12").Value = '  +11.51%  '
$ws.Range("D13").Value = '''0.0000275'
$ws.Range("E13").Value = '  +6.99%  '
$ws.Range("D14").Value = '''9.46'
$ws.Range("E14").Value = '  +4.79%  '
$ws.Range("D15").Value = '4.020.97'
$ws.Range("E15").Value = '  +5.85%  '
$ws.Range("D16").Value = '3.467.54'
$ws.Range("E16").Value = '  +5.96%  '
$ws.Range("D17").Value = '67.525.61'
$ws.Range("E17").Value = '  +5.80%  '
$ws.Range("E18").Value = '  +4.03%  '
$ws.Range("D19").Value = '''18.32'
$ws.Range("E19").Value = '  +5.85%  '
$ws.Range("D20").Value = '''11.83'
$ws.Range("E20").Value = '  +7.32%  '
$ws.Range("E21").Value = '  +6.17%  '
$ws.Range("D22").Value = '''405.89'
$ws.Range("E22").Value = '  +9.14%  '
$ws.Range("D23").Value = '''11.96'
$ws.Range("E23").Value = '  +9.05%  '
$ws.Range("D24").Value = '''3.91'
$ws.Range("E24").Value = '  +4.70%  '
$ws.Range("D25").Value = '''84.79'
$ws.Range("E25").Value = '  +5.94%  '
$ws.Range("E26").Value = '  +8.93%  '
$ws.Range("D27").Value = '''2.93'
$ws.Range("E27").Value = '  +10.80%  '
$ws.Range("E28").Value = '  +3.01%  '
$ws.Range("D29").Value = '''11.80'
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("D30").Value = '''8.66'
$ws.Range("E30").Value = '  +4.65%  '
$ws.Range("D31").Value = '''30.27'
$ws.Range("E31").Value = '  +5.71%  '
$ws.Range("D32").Value = '''679.57'
$ws.Range("E32").Value = '  +6.30%  '
$ws.Range("E33").Value = '  +3.86%  '
$ws.Range("D34").Value = '''11.70'
$ws.Range("E34").Value = '  +4.55%  '
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("D36").Value = '''58.94'
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = '0.0₃0831'
$ws.Range("E37").Value = '  +19.51%  '
$ws.Range("D38").Value = '''38.63'
$ws.Range("E38").Value = '  +6.15%  '
$ws.Range("D39").Value = '''0.404'
$ws.Range("E39").Value = '  +5.33%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '''3.45'
$ws.Range("E41").Value = '  +25.63%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.134'
$ws.Range("E42").Value = '  +8.20%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '''2.78'
$ws.Range("E43").Value = '  +14.20%  '
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '3.049.85'
$ws.Range("E45").Value = '  +5.03%  '
$ws.Range("E46").Value = '  +11.76%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '''3.32'
$ws.Range("E47").Value = '  +8.98%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0422'
$ws.Range("E48").Value = '  +7.26%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '''2.74'
$ws.Range("E49").Value = '  +4.00%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Value = '''8.75'
$ws.Range("E50").Value = '  +12.42%  '
$ws.Range("E51").Value = '  +3.80%  '
